$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TXG")

# Row 6: "Change in inventories"
$ws.Range("B6").Value = 96000000.0
$ws.Range("C6").Value = 81000000.0
$ws.Range("D6").Value = 55020000.0
$ws.Range("E6").Value = 31851000.0
$ws.Range("F6").Value = 10265000.0

# Row 8: "Change in payables and accrued liability"
$ws.Range("B8").Value = 38000000.0
$ws.Range("C8").Value = 46000000.0
$ws.Range("D8").Value = 45301000.0
$ws.Range("E8").Value = 44434000.0
$ws.Range("F8").Value = 40889000.0
